$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.803.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.213.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.94%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.99"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +7.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0960"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.72%  "
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.541.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.866"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.213.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.638.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.91%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.36%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0737"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.62%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0309"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +23.01%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.203"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("B45").Value = "SynthetixNetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.102"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.73%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.24%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.07%  "
